# Edits three figure-caption text boxes that live inside the single big
# group shape ("组合 4") on slide 1.
#
#   Shape id 3078 "TextBox 1"  : "(a) DSI-11 Ex Vivo"        -> split into 3 runs
#   Shape id 3079 "TextBox 12" : "(b) DSI-11 b10k In Vivo"   -> hyphenate "b10k"
#   Shape id 27   "TextBox 12" : "(c) DSI-11 b7k In Vivo"    -> hyphenate "b7k"
#                                 and move it next to "DSI-11"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

# ---------------------------------------------------------------------
# (a) DSI-11 Ex Vivo  ->  "(a) " / "DSI-11 Ex " / "Vivo"   (3 runs)
# ---------------------------------------------------------------------
$shA = $g.GroupItems.Item(5)
$trA = $shA.TextFrame.TextRange
# Re-asserting Bold (already True) on each sub-range forces PowerPoint to
# split the run at these boundaries without altering any visible formatting.
$trA.Characters(1, 4).Font.Bold = $true     # "(a) "
$trA.Characters(5, 10).Font.Bold = $true    # "DSI-11 Ex "
$trA.Characters(15, 4).Font.Bold = $true    # "Vivo"

# ---------------------------------------------------------------------
# (b) DSI-11 b10k In Vivo  ->  (b) DSI-11-b10k In Vivo
# ---------------------------------------------------------------------
$shB = $g.GroupItems.Item(21)
$trB = $shB.TextFrame.TextRange
$trB.Characters(5, 12).Text = "DSI-11-b10k "

# ---------------------------------------------------------------------
# (c) DSI-11 b7k In Vivo  ->  (c) DSI-11-b7k In Vivo
# Shrink the later run first so the box never transiently grows to two
# lines (which would perturb the autofit height by a rounding EMU).
# ---------------------------------------------------------------------
$shC = $g.GroupItems.Item(23)
$trC = $shC.TextFrame.TextRange
$trC.Characters(12, 7).Text = "In "
$trC.Characters(5, 7).Text = "DSI-11-b7k "
